$wb = $excel.ActiveWorkbook
$before = $wb.Worksheets.Item("Progress Tracking")
$new = $wb.Worksheets.Add($before)
$new.Name = "user"

$new.Range("A1").Value = "userid"
$new.Range("B1").Value = "name"
$new.Range("A2").Value = 7843
$new.Range("B2").Value = "sajay"

$new.Range("A1:B1").Font.Bold = $true

$new.Activate() | Out-Null
$new.Range("J21").Select() | Out-Null
